$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.9371385197543067
$ws.Range("J2").Value = 0.9371385197543067
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05601
$ws.Range("N2").Value = 0.16803
$ws.Range("O2").Value = 0.02710547761971223
$ws.Range("P2").Value = 0.02710547761971223
$ws.Range("Q2").Value = 0.77274741664
$ws.Range("R2").Value = 6.954726749760001
$ws.Range("S2").Value = 0.02540158717377061
$ws.Range("T2").Value = 0.02540158717377061
$ws.Range("I3").Value = 0.9371385197543067
$ws.Range("J3").Value = 0.9371385197543067
$ws.Range("N3").Value = 5.594253
$ws.Range("O3").Value = 0.902427539668559
$ws.Range("P3").Value = 0.9024275396685592
$ws.Range("S3").Value = 0.8456996087105143
$ws.Range("T3").Value = 0.8456996087105144
$ws.Range("I4").Value = 0.9371385197543067
$ws.Range("J4").Value = 0.9371385197543067
$ws.Range("M4").Value = 0.145611
$ws.Range("N4").Value = 0.436833
$ws.Range("O4").Value = 0.07046698271172858
$ws.Range("P4").Value = 0.07046698271172858
$ws.Range("Q4").Value = 2.008936334304
$ws.Range("R4").Value = 18.080427008736
$ws.Range("S4").Value = 0.06603732387002165
$ws.Range("T4").Value = 0.06603732387002165
$ws.Range("G5").Value = 0.9254496666666667
$ws.Range("H5").Value = 2.776349
$ws.Range("I5").Value = 0.0628614802456932
$ws.Range("J5").Value = 0.06286148024569319
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05601
$ws.Range("N5").Value = 0.16803
$ws.Range("O5").Value = 0.02710547761971223
$ws.Range("P5").Value = 0.02710547761971223
$ws.Range("Q5").Value = 0.05183443583
$ws.Range("R5").Value = 0.4665099224700001
$ws.Range("S5").Value = 0.001703890445941619
$ws.Range("T5").Value = 0.001703890445941619
$ws.Range("G6").Value = 0.9254496666666667
$ws.Range("H6").Value = 2.776349
$ws.Range("I6").Value = 0.0628614802456932
$ws.Range("J6").Value = 0.06286148024569319
$ws.Range("N6").Value = 5.594253
$ws.Range("O6").Value = 0.902427539668559
$ws.Range("P6").Value = 0.9024275396685592
$ws.Range("Q6").Value = 1.725733191366333
$ws.Range("R6").Value = 15.531598722297
$ws.Range("S6").Value = 0.05672793095804464
$ws.Range("T6").Value = 0.05672793095804463
$ws.Range("G7").Value = 0.9254496666666667
$ws.Range("H7").Value = 2.776349
$ws.Range("I7").Value = 0.0628614802456932
$ws.Range("J7").Value = 0.06286148024569319
$ws.Range("M7").Value = 0.145611
$ws.Range("N7").Value = 0.436833
$ws.Range("O7").Value = 0.07046698271172858
$ws.Range("P7").Value = 0.07046698271172858
$ws.Range("Q7").Value = 0.134755651413
$ws.Range("R7").Value = 1.212800862717
$ws.Range("S7").Value = 0.00442965884170693
$ws.Range("T7").Value = 0.004429658841706929
